# This edit swaps the two message blocks in the sheet:
#   - Block A: rows 2-15  (message_id = 0, "Critical low battery ..." sentence, 14 tokens)
#   - Block B: rows 16-22 (message_id = 1, "Compass Error ..." sentence, 7 tokens)
# After the edit, Block B comes first (rows 2-8, message_id reset to 0) followed
# by Block A (rows 9-22, message_id reset to 1). token_index values travel with
# their original token/row, only the message_id column is renumbered based on
# the block's new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 12   # columns A..L

# Capture the two blocks (including all columns A-L) before overwriting anything.
$blockA = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(15, $lastCol)).Value2
$blockB = $ws.Range($ws.Cells.Item(16, 1), $ws.Cells.Item(22, $lastCol)).Value2

$rowsA = 14
$rowsB = 7

# Write Block B first, at rows 2..8, with message_id forced to 0
for ($i = 1; $i -le $rowsB; $i++) {
    $destRow = 1 + $i
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $blockB[$i, $c]
        if ($c -eq 1) {
            $val = 0
        }
        $ws.Cells.Item($destRow, $c).Value = $val
    }
}

# Then write Block A at rows 9..22, with message_id forced to 1
for ($i = 1; $i -le $rowsA; $i++) {
    $destRow = 8 + $i
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $blockA[$i, $c]
        if ($c -eq 1) {
            $val = 1
        }
        $ws.Cells.Item($destRow, $c).Value = $val
    }
}
